$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.029.76'
$ws.Range('E2').Value = '  -7.14%  '
$ws.Range('D3').Value = '1.416.06'
$ws.Range('E3').Value = '  -7.55%  '
$ws.Range('D4').Formula = "'1.000"
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Formula = "'1.000"
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Formula = "'275.34"
$ws.Range('E6').Value = '  -4.74%  '
$ws.Range('D7').Formula = "'0.3683"
$ws.Range('E7').Value = '  -5.37%  '
$ws.Range('D8').Formula = "'0.3125"
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').Formula = "'39.71"
$ws.Range('E9').Value = '  -7.18%  '
$ws.Range('D10').Formula = "'1.038"
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('D11').Formula = "'0.06513"
$ws.Range('E11').Value = '  -8.87%  '
$ws.Range('D12').Formula = "'1.000"
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Formula = "'5.491"
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('D14').Formula = "'17.72"
$ws.Range('E14').Value = '  -2.58%  '
$ws.Range('D15').Formula = "'6.205"
$ws.Range('E15').Value = '  -5.47%  '
$ws.Range('D16').Value = '1.416.30'
$ws.Range('E16').Value = '  -7.75%  '
$ws.Range('D17').Formula = "'0.00001021"
$ws.Range('E17').Value = '  -5.91%  '
$ws.Range('D18').Formula = "'0.05707"
$ws.Range('E18').Value = '  -13.82%  '
$ws.Range('D19').Formula = "'0.9997"
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').Formula = "'71.06"
$ws.Range('E20').Value = '  -15.32%  '
$ws.Range('D21').Formula = "'5.624"
$ws.Range('E21').Value = '  -7.78%  '
$ws.Range('E22').Value = '  -3.91%  '
$ws.Range('D23').Formula = "'11.14"
$ws.Range('E23').Value = '  +3.68%  '
$ws.Range('D24').Formula = "'2.260"
$ws.Range('E24').Value = '  -5.02%  '
$ws.Range('D25').Value = '20.051.56'
$ws.Range('E25').Value = '  -7.01%  '
$ws.Range('D26').Formula = "'2.268"
$ws.Range('E26').Value = '  -3.69%  '
$ws.Range('D27').Formula = "'135.43"
$ws.Range('E27').Value = '  -10.08%  '
$ws.Range('D28').Formula = "'17.14"
$ws.Range('E28').Value = '  -6.75%  '
$ws.Range('D29').Value = '1.577.62'
$ws.Range('E29').Value = '  -7.70%  '
$ws.Range('D30').Formula = "'109.85"
$ws.Range('E30').Value = '  -5.91%  '
$ws.Range('D31').Formula = "'3.967"
$ws.Range('E31').Value = '  -18.38%  '
$ws.Range('D32').Formula = "'5.366"
$ws.Range('E32').Value = '  -10.79%  '
$ws.Range('D33').Formula = "'0.8292"
$ws.Range('E33').Value = '  -12.66%  '
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('D35').Formula = "'8.481"
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Formula = "'1.476"
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('D37').Formula = "'0.05913"
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').Formula = "'4.884"
$ws.Range('E38').Value = '  -5.03%  '
$ws.Range('D39').Formula = "'0.9994"
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Formula = "'0.02073"
$ws.Range('E40').Value = '  -5.94%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Formula = "'10.61"
$ws.Range('E41').Value = '  -5.38%  '
$ws.Range('D42').Formula = "'0.1916"
$ws.Range('E42').Value = '  -5.26%  '
$ws.Range('D43').Formula = "'1.096"
$ws.Range('E43').Value = '  -6.89%  '
$ws.Range('D44').Formula = "'0.5311"
$ws.Range('E44').Value = '  -7.63%  '
$ws.Range('D45').Formula = "'12.32"
$ws.Range('E45').Value = '  -6.08%  '
$ws.Range('D46').Formula = "'3.534"
$ws.Range('E46').Value = '  -4.95%  '
$ws.Range('D47').Formula = "'0.5171"
$ws.Range('E47').Value = '  -6.56%  '
$ws.Range('D48').Formula = "'114.40"
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('D49').Formula = "'1.772"
$ws.Range('E49').Value = '  -5.91%  '
$ws.Range('D50').Formula = "'1.043"
$ws.Range('E50').Value = '  -10.12%  '
$ws.Range('D51').Formula = "'1.000"
$ws.Range('E51').Value = '  -0.27%  '
